# Applies the "Automatic update of files." commit:
#  - bumps the "Förändrad" (column C) date on every data row from 46066 to 46070
#  - re-sorts the data rows (4..113) into a new order (content for each
#    "Beteckning" moves as a whole row, formulas included)
#
# Strategy: bulk-read the plain-value block (A:R) with .Value2 (full double
# precision, no round-trip loss) and the formula block (S:Z) with .Formula
# (uniformly returns literal values or formula text), permute the rows in
# memory according to the mapping below, patch column C, then bulk-write
# both blocks back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 113

# For each new row (index 0 => row 2, index 1 => row 3, ...), the number of
# the OLD row whose content should be placed there.
$oldRowForNewRow = @(
  2,3,5,4,6,13,7,8,9,10,11,12,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,
  110,111,112,113,32,31,33,35,37,38,50,51,52,39,41,40,42,44,45,43,59,47,46,48,49,
  53,54,55,71,56,36,57,77,58,82,84,34,78,60,61,79,62,63,64,81,65,66,67,68,69,70,
  72,73,74,75,76,80,83,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,101,102,
  103,104,105,106,107,108,109
)

# --- Read the two blocks as 2D arrays, 1-based: (row, col) ---
$valRange = $ws.Range("A$firstDataRow`:R$lastDataRow")
$valArr = $valRange.Value2

$formRange = $ws.Range("S$firstDataRow`:Z$lastDataRow")
$formArr = $formRange.Formula

$numRows = $lastDataRow - $firstDataRow + 1
$numValCols = 18
$numFormCols = 8

# --- Build the permuted output arrays ---
# NOTE: arrays coming back from Range.Value2 / Range.Formula are 1-based
# (classic VBA SAFEARRAY style: index 1..N), but a freshly `New-Object`ed
# .NET array is 0-based. Read with 1-based indices, write with 0-based ones.
$newValArr = New-Object 'object[,]' $numRows,$numValCols
$newFormArr = New-Object 'object[,]' $numRows,$numFormCols

for ($i = 0; $i -lt $numRows; $i++) {
    $oldRow = $oldRowForNewRow[$i]
    $oldIdx = $oldRow - $firstDataRow + 1

    for ($c = 1; $c -le $numValCols; $c++) {
        $newValArr[$i,($c-1)] = $valArr[$oldIdx,$c]
    }
    for ($c = 1; $c -le $numFormCols; $c++) {
        $newFormArr[$i,($c-1)] = $formArr[$oldIdx,$c]
    }

    # Column C ("Förändrad") is the 3rd column of the A:R block -> bump to 46070
    $newValArr[$i,2] = 46070
}

# --- Write the reordered blocks back ---
$ws.Range("A$firstDataRow`:R$lastDataRow").Value2 = $newValArr
$ws.Range("S$firstDataRow`:Z$lastDataRow").Formula = $newFormArr
